# Update Section_A and Section_B timetable sheets with revised course
# scheduling slots.

$wb = $excel.ActiveWorkbook

# ---- Section_A ----
$wsA = $wb.Worksheets.Item("Section_A")

$wsA.Range("B2").Value = "DS302"
$wsA.Range("C2").Value = "DS302"
$wsA.Range("D2").Value = "Free"
$wsA.Range("E2").Value = "Free"

$wsA.Range("C3").Value = "DS303"
$wsA.Range("D3").Value = "CS307"
$wsA.Range("E3").Value = "Free"
$wsA.Range("F3").Value = "Free"

$wsA.Range("B5").Value = "Free"
$wsA.Range("C5").Value = "Free"
$wsA.Range("D5").Value = "DS302"
$wsA.Range("E5").Value = "DS303"

$wsA.Range("B6").Value = "DS303 (Tutorial)"
$wsA.Range("C6").Value = "Free"
$wsA.Range("D6").Value = "DS302 (Tutorial)"

$wsA.Range("F8").Value = "Free"

# ---- Section_B ----
$wsB = $wb.Worksheets.Item("Section_B")

$wsB.Range("B2").Value = "DS302"
$wsB.Range("C2").Value = "Free"
$wsB.Range("D2").Value = "Free"
$wsB.Range("E2").Value = "Free"
$wsB.Range("F2").Value = "CS307"

$wsB.Range("C3").Value = "Free"
$wsB.Range("E3").Value = "CS307"

$wsB.Range("C5").Value = "Free"
$wsB.Range("D5").Value = "DS303"
$wsB.Range("E5").Value = "DS303"
$wsB.Range("F5").Value = "DS302"

$wsB.Range("D6").Value = "Free"
$wsB.Range("F6").Value = "Free"

$wsB.Range("B7").Value = "CS307"
$wsB.Range("C7").Value = "Free"
$wsB.Range("D7").Value = "DS302"
$wsB.Range("E7").Value = "Free"
$wsB.Range("F7").Value = "DS303"

$wsB.Range("C8").Value = "DS303 (Tutorial)"
$wsB.Range("D8").Value = "DS302 (Tutorial)"

$wb.Save()
